# Refresh the crypto price/volume snapshot (GitHub Actions scheduled update).
# Price/Volume(1h) columns are stored as literal text (e.g. "30.085.77"), so
# values are written with a leading apostrophe where they'd otherwise be
# auto-parsed as numbers, then ClearFormats() strips the @ text format Excel
# stamps on forced-text numeric-looking entries, keeping the original (no
# explicit style) cell formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "30.085.77"
    "D3" = "1.876.90"
    "E3" = "  -2.27%  "
    "E4" = "  +0.32%  "
    "D5" = "'319.68"
    "E5" = "  -3.15%  "
    "E6" = "  +0.23%  "
    "D7" = "'0.5038"
    "E7" = "  -3.46%  "
    "D8" = "'0.3956"
    "E8" = "  -3.26%  "
    "D9" = "'0.08211"
    "E9" = "  -3.93%  "
    "D10" = "'42.03"
    "E10" = "  -2.43%  "
    "D11" = "'1.092"
    "E11" = "  -3.19%  "
    "D12" = "'23.59"
    "E12" = "  +5.45%  "
    "D13" = "1.874.60"
    "E13" = "  -2.04%  "
    "D14" = "'6.293"
    "E14" = "  -1.96%  "
    "D15" = "'7.198"
    "E15" = "  -2.90%  "
    "E16" = "  +0.31%  "
    "D17" = "'91.73"
    "E17" = "  -3.87%  "
    "E18" = "  -2.66%  "
    "D19" = "'0.06462"
    "E19" = "  -3.40%  "
    "D20" = "'18.07"
    "E20" = "  -1.92%  "
    "E21" = "  +0.28%  "
    "D22" = "30.084.29"
    "D23" = "'5.835"
    "E23" = "  -2.97%  "
    "D24" = "'11.14"
    "E24" = "  -2.30%  "
    "D25" = "'2.155"
    "E25" = "  -2.54%  "
    "D26" = "2.094.70"
    "E26" = "  -1.93%  "
    "D27" = "'161.11"
    "E27" = "  +0.79%  "
    "E28" = "  -0.10%  "
    "D29" = "'2.243"
    "E29" = "  -8.37%  "
    "D30" = "'127.61"
    "E30" = "  -1.27%  "
    "D31" = "'1.073"
    "E31" = "  -0.96%  "
    "E32" = "  -2.61%  "
    "D33" = "'5.916"
    "E33" = "  -2.26%  "
    "D34" = "'3.699"
    "E34" = "  +1.76%  "
    "D35" = "'0.02425"
    "E35" = "  -2.86%  "
    "D36" = "'5.267"
    "E36" = "  +1.64%  "
    "D37" = "'0.06346"
    "D38" = "'0.2130"
    "E38" = "  -3.57%  "
    "E39" = "  -5.08%  "
    "D40" = "'8.497"
    "E40" = "  -4.79%  "
    "D41" = "'1.217"
    "E41" = "  -2.87%  "
    "D42" = "'0.6283"
    "E42" = "  -4.14%  "
    "D43" = "'11.26"
    "E44" = "  +0.28%  "
    "D45" = "'13.25"
    "E45" = "  -0.27%  "
    "D46" = "'0.5904"
    "E46" = "  -4.18%  "
    "D47" = "'2.083"
    "E47" = "  +0.11%  "
    "D48" = "'3.629"
    "E48" = "  -3.57%  "
    "B49" = "EOS"
    "C49" = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
    "D49" = "'1.209"
    "E49" = "  -3.28%  "
    "B50" = "Quant"
    "C50" = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
    "D50" = "'122.17"
    "E50" = "  -1.97%  "
    "D51" = "'77.35"
    "E51" = "  -3.12%  "
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.Value = $updates[$ref]
    $cell.ClearFormats()
}
